# Auto-generated edit script applying the Marilith_Profits value updates
$wb = $excel.ActiveWorkbook

$changes = @(
    @{ Sheet = "ALC"; Row = 2; Col = 8; Value = 1278 },
    @{ Sheet = "ALC"; Row = 2; Col = 9; Value = 994 },
    @{ Sheet = "ALC"; Row = 2; Col = 11; Value = 994 },
    @{ Sheet = "ALC"; Row = 2; Col = 13; Value = -881 },
    @{ Sheet = "ALC"; Row = 6; Col = 8; Value = 307.625 },
    @{ Sheet = "ALC"; Row = 6; Col = 9; Value = 284.7 },
    @{ Sheet = "ALC"; Row = 6; Col = 10; Value = 345.83334 },
    @{ Sheet = "ALC"; Row = 6; Col = 11; Value = 854.0999999999999 },
    @{ Sheet = "ALC"; Row = 6; Col = 12; Value = 1037.50002 },
    @{ Sheet = "ALC"; Row = 6; Col = 13; Value = -742.0999999999999 },
    @{ Sheet = "ALC"; Row = 6; Col = 14; Value = -1261.50002 },
    @{ Sheet = "ALC"; Row = 15; Col = 8; Value = 1006.5455 },
    @{ Sheet = "ALC"; Row = 15; Col = 9; Value = 1006.5455 },
    @{ Sheet = "ALC"; Row = 15; Col = 11; Value = 3019.6365 },
    @{ Sheet = "ALC"; Row = 15; Col = 13; Value = -2850.6365 },
    @{ Sheet = "ALC"; Row = 40; Col = 8; Value = 4386.625 },
    @{ Sheet = "ALC"; Row = 40; Col = 9; Value = 3500 },
    @{ Sheet = "ALC"; Row = 40; Col = 11; Value = 3500 },
    @{ Sheet = "ALC"; Row = 40; Col = 13; Value = -3325 },
    @{ Sheet = "ALC"; Row = 51; Col = 8; Value = 4997.5 },
    @{ Sheet = "ALC"; Row = 51; Col = 10; Value = 4998.5 },
    @{ Sheet = "ALC"; Row = 51; Col = 12; Value = 4998.5 },
    @{ Sheet = "ALC"; Row = 51; Col = 14; Value = -5966.5 },
    @{ Sheet = "ALC"; Row = 62; Col = 8; Value = 5432.6665 },
    @{ Sheet = "ALC"; Row = 62; Col = 9; Value = 5112.0625 },
    @{ Sheet = "ALC"; Row = 62; Col = 11; Value = 5112.0625 },
    @{ Sheet = "ALC"; Row = 62; Col = 13; Value = -4488.0625 },
    @{ Sheet = "ALC"; Row = 65; Col = 8; Value = 5432.6665 },
    @{ Sheet = "ALC"; Row = 65; Col = 9; Value = 5112.0625 },
    @{ Sheet = "ALC"; Row = 65; Col = 11; Value = 25560.3125 },
    @{ Sheet = "ALC"; Row = 65; Col = 13; Value = -22440.3125 },
    @{ Sheet = "ALC"; Row = 86; Col = 8; Value = 8684.6875 },
    @{ Sheet = "ALC"; Row = 86; Col = 9; Value = 3000 },
    @{ Sheet = "ALC"; Row = 86; Col = 10; Value = 18159.166 },
    @{ Sheet = "ALC"; Row = 86; Col = 11; Value = 3000 },
    @{ Sheet = "ALC"; Row = 86; Col = 12; Value = 18159.166 },
    @{ Sheet = "ALC"; Row = 86; Col = 13; Value = -1877 },
    @{ Sheet = "ALC"; Row = 86; Col = 14; Value = -20405.166 },
    @{ Sheet = "ALC"; Row = 89; Col = 8; Value = 8684.6875 },
    @{ Sheet = "ALC"; Row = 89; Col = 9; Value = 3000 },
    @{ Sheet = "ALC"; Row = 89; Col = 10; Value = 18159.166 },
    @{ Sheet = "ALC"; Row = 89; Col = 11; Value = 15000 },
    @{ Sheet = "ALC"; Row = 89; Col = 12; Value = 90795.83 },
    @{ Sheet = "ALC"; Row = 89; Col = 13; Value = -9384 },
    @{ Sheet = "ALC"; Row = 89; Col = 14; Value = -102027.83 },
    @{ Sheet = "ALC"; Row = 92; Col = 8; Value = 786.55554 },
    @{ Sheet = "ALC"; Row = 92; Col = 9; Value = 497.14285 },
    @{ Sheet = "ALC"; Row = 92; Col = 11; Value = 497.14285 },
    @{ Sheet = "ALC"; Row = 92; Col = 13; Value = 750.85715 },
    @{ Sheet = "ALC"; Row = 100; Col = 8; Value = 2454.7144 },
    @{ Sheet = "ALC"; Row = 100; Col = 9; Value = 2636.8 },
    @{ Sheet = "ALC"; Row = 100; Col = 10; Value = 1999.5 },
    @{ Sheet = "ALC"; Row = 100; Col = 11; Value = 2636.8 },
    @{ Sheet = "ALC"; Row = 100; Col = 12; Value = 1999.5 },
    @{ Sheet = "ALC"; Row = 100; Col = 13; Value = -2095.8 },
    @{ Sheet = "ALC"; Row = 100; Col = 14; Value = -3081.5 },
    @{ Sheet = "ALC"; Row = 106; Col = 8; Value = 4199.6665 },
    @{ Sheet = "ALC"; Row = 106; Col = 9; Value = 4299.5 },
    @{ Sheet = "ALC"; Row = 106; Col = 10; Value = 4000 },
    @{ Sheet = "ALC"; Row = 106; Col = 11; Value = 4299.5 },
    @{ Sheet = "ALC"; Row = 106; Col = 12; Value = 4000 },
    @{ Sheet = "ALC"; Row = 106; Col = 13; Value = -3668.5 },
    @{ Sheet = "ALC"; Row = 106; Col = 14; Value = -5262 },
    @{ Sheet = "ALC"; Row = 125; Col = 8; Value = 2302.25 },
    @{ Sheet = "ALC"; Row = 125; Col = 9; Value = 916.8570999999999 },
    @{ Sheet = "ALC"; Row = 125; Col = 11; Value = 8251.713899999999 },
    @{ Sheet = "ALC"; Row = 125; Col = 13; Value = -5791.713899999999 },
    @{ Sheet = "ALC"; Row = 137; Col = 8; Value = 2069.5264 },
    @{ Sheet = "ALC"; Row = 137; Col = 9; Value = 1403.3572 },
    @{ Sheet = "ALC"; Row = 137; Col = 10; Value = 3934.8 },
    @{ Sheet = "ALC"; Row = 137; Col = 11; Value = 4210.071599999999 },
    @{ Sheet = "ALC"; Row = 137; Col = 12; Value = 11804.4 },
    @{ Sheet = "ALC"; Row = 137; Col = 13; Value = -1660.071599999999 },
    @{ Sheet = "ALC"; Row = 137; Col = 14; Value = -16904.4 },
    @{ Sheet = "ARM"; Row = 32; Col = 8; Value = 9826.525 },
    @{ Sheet = "ARM"; Row = 32; Col = 9; Value = 9826.525 },
    @{ Sheet = "ARM"; Row = 32; Col = 11; Value = 9826.525 },
    @{ Sheet = "ARM"; Row = 32; Col = 13; Value = -9539.525 },
    @{ Sheet = "ARM"; Row = 45; Col = 8; Value = 0 },
    @{ Sheet = "ARM"; Row = 45; Col = 9; Value = 0 },
    @{ Sheet = "ARM"; Row = 45; Col = 11; Value = 0 },
    @{ Sheet = "ARM"; Row = 45; Col = 13; Clear = $true },
    @{ Sheet = "ARM"; Row = 61; Col = 8; Value = 1395.125 },
    @{ Sheet = "ARM"; Row = 61; Col = 9; Value = 1395.125 },
    @{ Sheet = "ARM"; Row = 61; Col = 11; Value = 1395.125 },
    @{ Sheet = "ARM"; Row = 61; Col = 13; Value = -1183.125 },
    @{ Sheet = "ARM"; Row = 97; Col = 8; Value = 1556.5555 },
    @{ Sheet = "ARM"; Row = 97; Col = 9; Value = 853.3333 },
    @{ Sheet = "ARM"; Row = 97; Col = 10; Value = 2963 },
    @{ Sheet = "ARM"; Row = 97; Col = 11; Value = 853.3333 },
    @{ Sheet = "ARM"; Row = 97; Col = 12; Value = 2963 },
    @{ Sheet = "ARM"; Row = 97; Col = 13; Value = -357.3333 },
    @{ Sheet = "ARM"; Row = 97; Col = 14; Value = -3955 },
    @{ Sheet = "ARM"; Row = 136; Col = 8; Value = 1395.125 },
    @{ Sheet = "ARM"; Row = 136; Col = 9; Value = 1395.125 },
    @{ Sheet = "ARM"; Row = 136; Col = 11; Value = 4185.375 },
    @{ Sheet = "ARM"; Row = 136; Col = 13; Value = -1635.375 },
    @{ Sheet = "BSM"; Row = 75; Col = 8; Value = 24026.857 },
    @{ Sheet = "BSM"; Row = 75; Col = 9; Value = 8031.3335 },
    @{ Sheet = "BSM"; Row = 75; Col = 11; Value = 8031.3335 },
    @{ Sheet = "BSM"; Row = 75; Col = 13; Value = -7095.3335 },
    @{ Sheet = "BSM"; Row = 78; Col = 8; Value = 24026.857 },
    @{ Sheet = "BSM"; Row = 78; Col = 9; Value = 8031.3335 },
    @{ Sheet = "BSM"; Row = 78; Col = 11; Value = 24094.0005 },
    @{ Sheet = "BSM"; Row = 78; Col = 13; Value = -19414.0005 },
    @{ Sheet = "BSM"; Row = 107; Col = 8; Value = 2271.6 },
    @{ Sheet = "BSM"; Row = 107; Col = 9; Value = 1119.3334 },
    @{ Sheet = "BSM"; Row = 107; Col = 11; Value = 1119.3334 },
    @{ Sheet = "BSM"; Row = 107; Col = 13; Value = 800.6666 },
    @{ Sheet = "BSM"; Row = 124; Col = 8; Value = 79797 },
    @{ Sheet = "BSM"; Row = 124; Col = 10; Value = 79797 },
    @{ Sheet = "BSM"; Row = 124; Col = 12; Value = 79797 },
    @{ Sheet = "BSM"; Row = 124; Col = 14; Value = -89617 },
    @{ Sheet = "BSM"; Row = 129; Col = 8; Value = 0 },
    @{ Sheet = "BSM"; Row = 129; Col = 10; Value = 0 },
    @{ Sheet = "BSM"; Row = 129; Col = 12; Value = 0 },
    @{ Sheet = "BSM"; Row = 129; Col = 14; Clear = $true },
    @{ Sheet = "BSM"; Row = 134; Col = 8; Value = 11224.763 },
    @{ Sheet = "BSM"; Row = 134; Col = 9; Value = 8262.1 },
    @{ Sheet = "BSM"; Row = 134; Col = 10; Value = 14516.611 },
    @{ Sheet = "BSM"; Row = 134; Col = 11; Value = 24786.3 },
    @{ Sheet = "BSM"; Row = 134; Col = 12; Value = 43549.833 },
    @{ Sheet = "BSM"; Row = 134; Col = 13; Value = -22251.3 },
    @{ Sheet = "BSM"; Row = 134; Col = 14; Value = -48619.833 },
    @{ Sheet = "CRP"; Row = 44; Col = 8; Value = 23354.334 },
    @{ Sheet = "CRP"; Row = 44; Col = 9; Value = 23354.334 },
    @{ Sheet = "CRP"; Row = 44; Col = 10; Value = 0 },
    @{ Sheet = "CRP"; Row = 44; Col = 11; Value = 23354.334 },
    @{ Sheet = "CRP"; Row = 44; Col = 12; Value = 0 },
    @{ Sheet = "CRP"; Row = 44; Col = 13; Value = -22912.334 },
    @{ Sheet = "CRP"; Row = 44; Col = 14; Clear = $true },
    @{ Sheet = "CRP"; Row = 120; Col = 8; Value = 19999.334 },
    @{ Sheet = "CRP"; Row = 132; Col = 8; Value = 1704.9 },
    @{ Sheet = "CRP"; Row = 132; Col = 9; Value = 1442.8572 },
    @{ Sheet = "CRP"; Row = 132; Col = 10; Value = 2316.3333 },
    @{ Sheet = "CRP"; Row = 132; Col = 11; Value = 4328.571599999999 },
    @{ Sheet = "CRP"; Row = 132; Col = 12; Value = 6948.999899999999 },
    @{ Sheet = "CRP"; Row = 132; Col = 13; Value = -1798.571599999999 },
    @{ Sheet = "CRP"; Row = 132; Col = 14; Value = -12008.9999 },
    @{ Sheet = "CRP"; Row = 134; Col = 8; Value = 3852.1765 },
    @{ Sheet = "CRP"; Row = 134; Col = 9; Value = 3899.4 },
    @{ Sheet = "CRP"; Row = 134; Col = 11; Value = 11698.2 },
    @{ Sheet = "CRP"; Row = 134; Col = 13; Value = -9163.200000000001 },
    @{ Sheet = "CUL"; Row = 9; Col = 8; Value = 997 },
    @{ Sheet = "CUL"; Row = 9; Col = 10; Value = 997 },
    @{ Sheet = "CUL"; Row = 9; Col = 12; Value = 2991 },
    @{ Sheet = "CUL"; Row = 9; Col = 14; Value = -3439 },
    @{ Sheet = "CUL"; Row = 11; Col = 8; Value = 1313.8572 },
    @{ Sheet = "CUL"; Row = 11; Col = 10; Value = 2665.6667 },
    @{ Sheet = "CUL"; Row = 11; Col = 12; Value = 7997.000100000001 },
    @{ Sheet = "CUL"; Row = 11; Col = 14; Value = -8277.000100000001 },
    @{ Sheet = "CUL"; Row = 12; Col = 8; Value = 112.2 },
    @{ Sheet = "CUL"; Row = 12; Col = 9; Value = 69.5 },
    @{ Sheet = "CUL"; Row = 12; Col = 10; Value = 140.66667 },
    @{ Sheet = "CUL"; Row = 12; Col = 11; Value = 208.5 },
    @{ Sheet = "CUL"; Row = 12; Col = 12; Value = 422.00001 },
    @{ Sheet = "CUL"; Row = 12; Col = 13; Value = -35.5 },
    @{ Sheet = "CUL"; Row = 12; Col = 14; Value = -768.00001 },
    @{ Sheet = "GSM"; Row = 33; Col = 8; Value = 24999.5 },
    @{ Sheet = "GSM"; Row = 33; Col = 9; Value = 0 },
    @{ Sheet = "GSM"; Row = 33; Col = 10; Value = 24999.5 },
    @{ Sheet = "GSM"; Row = 33; Col = 11; Value = 0 },
    @{ Sheet = "GSM"; Row = 33; Col = 12; Value = 24999.5 },
    @{ Sheet = "GSM"; Row = 33; Col = 13; Clear = $true },
    @{ Sheet = "GSM"; Row = 33; Col = 14; Value = -25503.5 },
    @{ Sheet = "GSM"; Row = 102; Col = 8; Value = 3629.5 },
    @{ Sheet = "GSM"; Row = 102; Col = 10; Value = 0 },
    @{ Sheet = "GSM"; Row = 102; Col = 12; Value = 0 },
    @{ Sheet = "GSM"; Row = 102; Col = 14; Clear = $true },
    @{ Sheet = "GSM"; Row = 126; Col = 8; Value = 0 },
    @{ Sheet = "GSM"; Row = 126; Col = 9; Value = 0 },
    @{ Sheet = "GSM"; Row = 126; Col = 11; Value = 0 },
    @{ Sheet = "GSM"; Row = 126; Col = 13; Clear = $true },
    @{ Sheet = "LTW"; Row = 22; Col = 8; Value = 2134.6365 },
    @{ Sheet = "LTW"; Row = 22; Col = 9; Value = 1000 },
    @{ Sheet = "LTW"; Row = 22; Col = 10; Value = 2248.1 },
    @{ Sheet = "LTW"; Row = 22; Col = 11; Value = 1000 },
    @{ Sheet = "LTW"; Row = 22; Col = 12; Value = 2248.1 },
    @{ Sheet = "LTW"; Row = 22; Col = 13; Value = -705 },
    @{ Sheet = "LTW"; Row = 22; Col = 14; Value = -2838.1 },
    @{ Sheet = "LTW"; Row = 27; Col = 8; Value = 2134.6365 },
    @{ Sheet = "LTW"; Row = 27; Col = 9; Value = 1000 },
    @{ Sheet = "LTW"; Row = 27; Col = 10; Value = 2248.1 },
    @{ Sheet = "LTW"; Row = 27; Col = 11; Value = 1000 },
    @{ Sheet = "LTW"; Row = 27; Col = 12; Value = 2248.1 },
    @{ Sheet = "LTW"; Row = 27; Col = 13; Value = -893 },
    @{ Sheet = "LTW"; Row = 27; Col = 14; Value = -2462.1 },
    @{ Sheet = "LTW"; Row = 39; Col = 8; Value = 0 },
    @{ Sheet = "LTW"; Row = 39; Col = 10; Value = 0 },
    @{ Sheet = "LTW"; Row = 39; Col = 12; Value = 0 },
    @{ Sheet = "LTW"; Row = 39; Col = 14; Clear = $true },
    @{ Sheet = "LTW"; Row = 132; Col = 8; Value = 24196.334 },
    @{ Sheet = "LTW"; Row = 132; Col = 9; Value = 24772 },
    @{ Sheet = "LTW"; Row = 132; Col = 11; Value = 74316 },
    @{ Sheet = "LTW"; Row = 132; Col = 13; Value = -71786 },
    @{ Sheet = "LTW"; Row = 133; Col = 8; Value = 70000 },
    @{ Sheet = "LTW"; Row = 133; Col = 10; Value = 70000 },
    @{ Sheet = "LTW"; Row = 133; Col = 12; Value = 70000 },
    @{ Sheet = "LTW"; Row = 133; Col = 14; Value = -75060 },
    @{ Sheet = "LTW"; Row = 136; Col = 8; Value = 3583.4 },
    @{ Sheet = "LTW"; Row = 136; Col = 9; Value = 3378 },
    @{ Sheet = "LTW"; Row = 136; Col = 11; Value = 10134 },
    @{ Sheet = "LTW"; Row = 136; Col = 13; Value = -7584 },
    @{ Sheet = "WVR"; Row = 17; Col = 8; Value = 933.3333 },
    @{ Sheet = "WVR"; Row = 17; Col = 9; Value = 300 },
    @{ Sheet = "WVR"; Row = 17; Col = 10; Value = 1250 },
    @{ Sheet = "WVR"; Row = 17; Col = 11; Value = 300 },
    @{ Sheet = "WVR"; Row = 17; Col = 12; Value = 1250 },
    @{ Sheet = "WVR"; Row = 17; Col = 13; Value = -128 },
    @{ Sheet = "WVR"; Row = 17; Col = 14; Value = -1594 },
    @{ Sheet = "WVR"; Row = 132; Col = 8; Value = 3701.5 },
    @{ Sheet = "WVR"; Row = 132; Col = 9; Value = 3701.5 },
    @{ Sheet = "WVR"; Row = 132; Col = 11; Value = 11104.5 },
    @{ Sheet = "WVR"; Row = 132; Col = 13; Value = -8574.5 }
)

foreach ($chg in $changes) {
    $ws = $wb.Worksheets.Item($chg.Sheet)
    $cell = $ws.Cells.Item($chg.Row, $chg.Col)
    if ($chg.ContainsKey("Clear")) {
        $cell.ClearContents()
    } else {
        $cell.Value = $chg.Value
    }
}

Write-Output "Applied $($changes.Count) cell updates"